$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "description"
$ws.Range("A2").Value = "Production"
$ws.Range("B2").Value = "Production of the flow"
$ws.Range("A3").Value = "Consumption"
$ws.Range("B3").Value = "Consumption of the flow"

$ws.Columns.Item(1).EntireColumn.AutoFit()

$ws.Range("E8").Select()
